$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value (cascades recalculation of B6, C6, E6 automatically)
$ws.Range("B3").Value = 0.1378

# Move the active selection from G4 to B4
$ws.Range("B4").Select()

# New header row (row 9) reusing existing labels plus two new ones
$ws.Range("G9").Value = "Desired_Flow (mL_per_Day)"
$ws.Range("H9").Value = "Microsteps_period"
$ws.Range("I9").Value = "Constant"
$ws.Range("J9").Value = "Verification_Microsteps_Period"

# New data/formula row (row 10)
$ws.Range("G10").Value = 0.1
$ws.Range("H10").Formula = "=(G3*H3*J3*K3/(1440*60*0.715*I3))"
$ws.Range("I10").Formula = "=(H3*J3*K3/(1440*60*0.715*I3))"
$ws.Range("J10").Formula = "=I10*G10"

# New formula cell in row 14
$ws.Range("G14").Formula = "=4.96/(1.25*20)"
